# issue #494: add XLSX file reading
#
# Rebuilds Feuil1's data (new header row + numbers/formula/quoted text) and
# adds a second, empty worksheet ("Feuil2") that becomes the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Feuil1: replace the old 2-column text table with the new layout ----
# Row 1: headers
$ws1.Range("A1").Value = "titre1"
$ws1.Range("B1").Value = "titre2"
$ws1.Range("C1").Value = "titre 3"

# Row 2: a quote-prefixed text "2", a plain number, and a text value
$ws1.Range("A2").Value = "'2"
$ws1.Range("B2").Value = 3
$ws1.Range("C2").Value = "toto"

# Row 3: A3 is cleared out entirely, B3 is a number formatted as Text,
# C3 is a formula
$ws1.Range("A3").ClearContents()
$ws1.Range("B3").Value = 4
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("C3").Formula = "=3.7+B2"

# --- Add the second worksheet, placed right after Feuil1 ---------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Feuil2"

# --- Restore selections on each sheet, Feuil2 last so it ends up active -
[void]$ws1.Range("E2").Select()
[void]$ws2.Range("H8").Select()
